$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 76
$ws1.Cells.Item(5, 6).Value = 1664
$ws1.Cells.Item(6, 6).Value = 3265
$ws1.Cells.Item(7, 6).Value = 824
$ws1.Cells.Item(8, 6).Value = 2058
$ws1.Cells.Item(9, 6).Value = 1972
$ws1.Cells.Item(10, 6).Value = 1017
$ws1.Cells.Item(13, 6).Value = 1611
$ws1.Cells.Item(14, 6).Value = 347
$ws1.Cells.Item(16, 6).Value = 14
$ws1.Cells.Item(18, 6).Value = 67
$ws1.Cells.Item(19, 6).Value = 1443
$ws1.Cells.Item(20, 6).Value = 529
$ws1.Cells.Item(21, 6).Value = 632
$ws1.Cells.Item(22, 6).Value = 322
$ws1.Cells.Item(23, 6).Value = 10736
$ws1.Cells.Item(24, 6).Value = 10756
$ws1.Cells.Item(25, 6).Value = 856
$ws1.Cells.Item(26, 6).Value = 662
$ws1.Cells.Item(27, 6).Value = 1840
$ws1.Cells.Item(28, 6).Value = 149
$ws1.Cells.Item(29, 6).Value = 440

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 33

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 62

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 76
$ws4.Cells.Item(3, 6).Value = 62
$ws4.Cells.Item(7, 6).Value = 1664
$ws4.Cells.Item(8, 6).Value = 3265
$ws4.Cells.Item(9, 6).Value = 825
$ws4.Cells.Item(10, 6).Value = 2058
$ws4.Cells.Item(11, 6).Value = 1972
$ws4.Cells.Item(12, 6).Value = 1017
$ws4.Cells.Item(15, 6).Value = 1611
$ws4.Cells.Item(16, 6).Value = 347
$ws4.Cells.Item(18, 6).Value = 14
$ws4.Cells.Item(22, 6).Value = 67
$ws4.Cells.Item(23, 6).Value = 1443
$ws4.Cells.Item(24, 6).Value = 529
$ws4.Cells.Item(25, 6).Value = 632
$ws4.Cells.Item(26, 6).Value = 322
$ws4.Cells.Item(27, 6).Value = 10736
$ws4.Cells.Item(28, 6).Value = 10756
$ws4.Cells.Item(29, 6).Value = 856
$ws4.Cells.Item(30, 6).Value = 662
$ws4.Cells.Item(31, 6).Value = 1840
$ws4.Cells.Item(33, 6).Value = 33
$ws4.Cells.Item(34, 6).Value = 149
$ws4.Cells.Item(35, 6).Value = 440
